$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells that get their old numeric "sort index" value replaced by a freshly
# generated UUID string (and, as a consequence of typing a brand-new value
# straight into the cell, lose their previously-applied explicit style).
$updates = @{
    "A1"  = "08f8c612-2700-406e-88df-eae964a98f67"
    "B2"  = "c2ddea44-b60b-49ad-b919-7e40ecddcdb9"
    "C3"  = "2f14b53d-1bcc-4a1b-97ca-d08cfedbc31c"
    "C4"  = "5b0a804a-d0c0-45c7-9ef0-bf5a5bba271a"
    "C5"  = "7580d708-5314-4212-a30d-81fad3738b59"
    "B6"  = "59c4f253-c906-4e42-a81d-75a735e299ef"
    "C7"  = "6e56c599-9fb4-47e6-b7cd-191e59634c7b"
    "C8"  = "ee3fad65-2876-45c4-9e91-7ca21405df92"
    "C9"  = "230daedc-6335-416f-95b4-ecfd952f60dc"
    "A10" = "7cbb7091-747e-4d5c-82e5-d7da589888e9"
    "B11" = "2cc28716-4861-4ca0-b78e-aa9194c28e1f"
    "C12" = "29a5e408-08c7-4f66-9432-972dd4d644b3"
    "C13" = "2cde7d42-04a6-42ad-8d2d-7ed704185105"
    "C14" = "f48aab6b-64e7-4aa8-8b0e-f70c4afe8cec"
    "B15" = "c27fa7a9-1170-4224-b73d-eb1b9e8c58a8"
    "C16" = "b27a63ba-0104-4471-992a-a6503a940cd0"
    "C17" = "3117ee28-6d3f-4586-a027-7addd15d0c42"
    "C18" = "f1c96124-e005-467e-b2f6-86be90a8eb03"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.Clear()
    $cell.Value = $updates[$addr]
}

# The trailing formatting-only row 1000 is removed from the sheet.
$ws.Rows.Item(1000).Delete()

# Selection moved from C23 to E22.
$ws.Range("E22").Select()
